$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("B1").Value = "VEGA CARDENAS, ANGELICA LOURDES - 08:15AM - 12:00PM"
$ws.Range("C1").Value = "SOTELO GONZALES, CAMILA SOFÍA - 01:30PM - 10:30PM"

# Row 2
$ws.Range("B2").Value = "MARTINEZ PAZ, ROCIO ESPERANZA - 08:45AM - 05:45PM"
$ws.Range("C2").Value = "SALAS VILLANUEVA, JAMILA DASHA - 05:45PM - 09:30PM"

# Row 3
$ws.Range("B3").Value = "POBLETE SAIRE, FIORELLA ESTHER - 09:15AM - 01:00PM"
$ws.Range("C3").Value = "TITO LAURA, NANCY FIORELLA - 02:00PM - 11:00PM"

# Row 4
$ws.Range("B4").Value = "YOVERA ROBLES, VICTOR EDUARDO - 06:30AM - 10:15AM"
$ws.Range("C4").Value = "MONTEZUMA DEJO, EVELYN BRUNELLA - 10:15AM - 02:00PM"
$ws.Range("D4").Value = "PARICELA TINEO, JAIME DANIEL - 02:00PM - 05:45PM"
$ws.Range("E4").Value = "VILCAPOMA CHILIN, JULISSA JAZMIN - 06:00PM - 09:45PM"

# Row 5
$ws.Range("B5").Value = "ALVITE CORNEJO, ANGIE LUCERO - 07:00AM - 10:45AM"
$ws.Range("C5").Value = "RUIZ SANTOS, CIELO CRISTHINA - 11:00AM - 02:45PM"
$ws.Range("D5").Value = "RIVERA CARREÑO, DIANA DESIRÉE - 03:30PM - 07:15PM"

# Row 6
$ws.Range("B6").Value = "HUAMAN HUAMANI, ALEXIS JAVIER - 08:30AM - 12:15PM"
$ws.Range("C6").Value = "CAPCHA YARANGO, DAVID - 02:00PM - 05:45PM"
$ws.Range("D6").Value = "ARIAS MACHACUAY, SADELITH SORAGGI - 05:45PM - 09:30PM"

# Row 7
$ws.Range("B7").Value = "QUISPE MONDRAGÓN, JUAN ALFONSO - 08:00AM - 11:45AM"
$ws.Range("C7").Value = "BARRIENTOS JERI, MILAGROS NICOL - 12:00PM - 03:30PM"
$ws.Range("D7").Value = "MUÑOZ SOTOMAYOR, MIRIAN RAQUEL - 03:30PM - 07:15PM"

# Row 8
$ws.Range("B8").Value = "BRICEÑO LUNA, JESSICA ARACELI - 11:00AM - 10:00PM"
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""

# Row 9
$ws.Range("B9").Value = "YANQUI BRAVO, MIRIAN LUZ - 09:00AM - 12:45PM"
$ws.Range("C9").Value = "YACILA GRANDEZ, RODRIGO ANDRE - 02:00PM - 05:45PM"
$ws.Range("D9").Value = "FLORES PAREDES, LOURDES - 06:00PM - 08:45PM"

# Row 10
$ws.Range("B10").Value = "AYALA MORA, CECILIA ROSARIO - 11:15AM - 03:00PM"
$ws.Range("C10").Value = "CHAVEZ ONOFRE, CAMILA GERALDINE - 06:00PM - 10:30PM"
$ws.Range("D10").Value = ""

# Row 11
$ws.Range("B11").Value = "MEZA PEREZ, JUAN CRISTOFER - 09:15AM - 01:00PM"
$ws.Range("C11").Value = "BONILLA SÁNCHEZ, RAÚL FERNANDO - 03:45PM - 07:30PM"

# Row 12
$ws.Range("B12").Value = "HUAYANAY VELASCO, ATHINA - 10:00AM - 01:45PM"
$ws.Range("C12").Value = "MENDOZA CRUZ, LILIANA LILIANA - 04:30PM - 08:15PM"

# Row 13
$ws.Range("B13").Value = "GOMEZ ALBINO, IDALIA GIMENA - 10:15AM - 02:00PM"
$ws.Range("C13").Value = "MEDINA MARCELO, NAOMI ARIADNA - 04:30PM - 08:15PM"

# Row 14
$ws.Range("B14").Value = "ILDEFONSO MOTTA, JHOSSEP ANGELO - 10:15AM - 02:00PM"
$ws.Range("C14").Value = "MENDOZA DIEGO, ZAIDA VANESSA - 05:15PM - 09:00PM"

# Row 15
$ws.Range("B15").Value = "QUIQUIA MALLQUI, CYNTHIA ANGELLINE - 11:15AM - 03:00PM"
$ws.Range("C15").Value = "SUAREZ JARA, YENNIFER YUSSARA - 05:30PM - 09:15PM"
